# Append a new data row (row 66) to each of the four sheets, mirroring the
# existing layout (time / B / C / D / E hex-strings / F / G / H / I counts).
#
# Column A keeps the existing "YYYY-MM-DD HH:MM:SS" custom date format so it
# reuses the same style as the rows above it; columns B-E stay text (hex byte
# lists); columns F-I are plain numbers.

$wb = $excel.ActiveWorkbook

$rowsData = @(
    @{
        Sheet = 1
        A = [double]"45852.49320601852"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x44"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 324
        I = 15
    },
    @{
        Sheet = 2
        A = [double]"45852.49320601852"
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x54"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 340
        I = 14
    },
    @{
        Sheet = 3
        A = [double]"45852.49320601852"
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x64"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 100
        I = 3
    },
    @{
        Sheet = 4
        A = [double]"45852.49320601852"
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x64"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 100
        I = 3
    }
)

foreach ($rd in $rowsData) {
    $ws = $wb.Worksheets.Item($rd.Sheet)

    $ws.Range("A66").Value2 = $rd.A
    $ws.Range("A66").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("B66").Value2 = $rd.B
    $ws.Range("C66").Value2 = $rd.C
    $ws.Range("D66").Value2 = $rd.D
    $ws.Range("E66").Value2 = $rd.E

    $ws.Range("F66").Value2 = $rd.F
    $ws.Range("G66").Value2 = $rd.G
    $ws.Range("H66").Value2 = $rd.H
    $ws.Range("I66").Value2 = $rd.I
}

Write-Output "row 66 appended to all sheets"
